# Update the Summary Table and Statistics
# - Header I1: "Folder" -> "Tools"
# - New WWAN/Windows test row (row 10) results are filled in:
#     F10 (Fail count), G10 (Fail count2) get real numbers
#     I10 test-item name corrected to airplane_wwan_windows
#     J10 date/time stamp added
#     K10 remark added
# - Move selection off the now-filled K10 cell back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value2 = "Tools"

$ws.Range("F10").Value2 = 603
$ws.Range("G10").Value2 = 397
$ws.Range("I10").Value2 = "airplane_wwan_windows"
$ws.Range("J10").Value2 = "20191106_041905"
$ws.Range("K10").Value2 = "WWAN module crash at the 604 times"

$ws.Range("A1").Select() | Out-Null
